# Realestate Update resale numbers 2023-06-15 22:31
# Appends a new data row (row 48) to the CityResaleNum sheet with the
# latest resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 48

# Columns A (Date) and D (Week) look like numbers/dates to Excel's
# auto-detection, so force them to text first (matching the source
# workbook, where these are stored as plain text strings), then drop the
# number-format override so no extra style is left behind on the cell.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2023-06-15"
$ws.Range("A$row").ClearFormats()

$ws.Range("B$row").Value = "22:30:50"

$ws.Range("C$row").Value = "Thursday"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "24"
$ws.Range("D$row").ClearFormats()

$ws.Range("E$row").Value = 121649
$ws.Range("F$row").Value = 132580
$ws.Range("G$row").Value = 161960
$ws.Range("H$row").Value = 132873
$ws.Range("I$row").Value = 176644
$ws.Range("J$row").Value = 114356
$ws.Range("K$row").Value = 200012
$ws.Range("L$row").Value = 224390
$ws.Range("M$row").Value = 174651
$ws.Range("N$row").Value = 102835
$ws.Range("O$row").Value = 39079
$ws.Range("P$row").Value = 34050
$ws.Range("Q$row").Value = 51694
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 36892
$ws.Range("T$row").Value = -1
